$wb = $excel.ActiveWorkbook

# "Generate Report for Handback" - refresh the handoff/handback timestamps
# for the most-recently processed file (row 2) on each language sheet.

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E2").Value = "2016-03-11 12:30:57"
$zhcn.Range("H2").Value = "2016-03-11 12:31:17"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E2").Value = "2016-03-11 12:31:00"
$dede.Range("H2").Value = "2016-03-11 12:31:25"
